# Fruta / hortaliza, semanal
# Insert a new daily price record at row 64 (pushing existing rows 64-123
# down to 65-124) for "Vega Central Mapocho de Santiago" / Frambuesa.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 64..123 down to 65..124, leaving a fresh blank row 64.
$ws.Rows.Item(64).Insert()

# Populate the new row 64 with the new record (same shape as the rows
# around it; only the date and volume differ from the old row 64).
$ws.Range("A64").Value = 9
$ws.Range("B64").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C64").Value = "Metropolitana"
$ws.Range("D64").Value = 44957
$ws.Range("E64").Value = 13
$ws.Range("F64").Value = "Fruta"
$ws.Range("G64").Value = 100101
$ws.Range("H64").Value = "Berries"
$ws.Range("I64").Value = 100101004
$ws.Range("J64").Value = "Frambuesa"
$ws.Range("K64").Value = "Sin especificar"
$ws.Range("L64").Value = "Primera"
$ws.Range("M64").Value = 400
$ws.Range("N64").Value = 7000
$ws.Range("O64").Value = 7000
$ws.Range("P64").Value = 7000
$ws.Range("Q64").Value = "$/bandeja 2 kilos"
$ws.Range("R64").Value = "Provincia de Curicó"
$ws.Range("S64").Value = 3500
$ws.Range("T64").Value = 2
